$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: I1 = "I0", J1 = "IF" (same style as the other header cells, e.g. H1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows 2-11: columns I and J ---
$data = @{
    2  = @(6, 7)
    3  = @(9, 9)
    4  = @(6, 6)
    5  = @(3, 3)
    6  = @(5, 6)
    7  = @(4, 4)
    8  = @(5, 5)
    9  = @(5, 5)
    10 = @(6, 6)
    11 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
